$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 80 (shifts old rows 80-120 down to 82-122)
$ws.Rows("80:81").Insert()

# New row 80
$ws.Range("A80").Value2 = 6
$ws.Range("B80").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C80").Value2 = "Metropolitana"
$ws.Range("D80").Value2 = 44510
$ws.Range("E80").Value2 = 13
$ws.Range("F80").Value2 = 100112001
$ws.Range("G80").Value2 = "Berenjena"
$ws.Range("H80").Value2 = "Sin especificar"
$ws.Range("I80").Value2 = "Primera"
$ws.Range("J80").Value2 = 220
$ws.Range("K80").Value2 = 7000
$ws.Range("L80").Value2 = 8000
$ws.Range("M80").Value2 = 7455
$ws.Range("N80").Value2 = "`$/caja 50 unidades"
$ws.Range("O80").Value2 = "Región de Arica y Parinacota"
$ws.Range("P80").Value2 = 149
$ws.Range("Q80").Value2 = 50
$ws.Range("R80").Value2 = "Hortaliza"

# New row 81
$ws.Range("A81").Value2 = 6
$ws.Range("B81").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C81").Value2 = "Metropolitana"
$ws.Range("D81").Value2 = 44510
$ws.Range("E81").Value2 = 13
$ws.Range("F81").Value2 = 100112001
$ws.Range("G81").Value2 = "Berenjena"
$ws.Range("H81").Value2 = "Sin especificar"
$ws.Range("I81").Value2 = "Primera"
$ws.Range("J81").Value2 = 200
$ws.Range("K81").Value2 = 12000
$ws.Range("L81").Value2 = 14000
$ws.Range("M81").Value2 = 13100
$ws.Range("N81").Value2 = "`$/caja 60 unidades"
$ws.Range("O81").Value2 = "Región de Arica y Parinacota"
$ws.Range("P81").Value2 = 218
$ws.Range("Q81").Value2 = 60
$ws.Range("R81").Value2 = "Hortaliza"
